$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value = "Sion Kim"

# Fix style of row 16 (E16:G16) to match the wrap-text/border style used by
# the rest of the "Expected Result" / "Method Inputs" / "Preconditions" columns
# (same style already present on E12:G15) before filling in values.
$ws.Range("E12:G12").Copy() | Out-Null
$ws.Range("E16:G16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 7 - __init__ / Attributes set to input values.
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "client_number: 1999`nfirst_name: ""Sion""`nlast_name: ""Kim""`nemail_address: ""skim14@rrc.ca"""
$ws.Range("G7").Value = "Object created with expected atttribute value based on method inputs."

# Expected Result column (G) for rows 8-11 filled first...
$ws.Range("G8").Value = "ValueError(""client_number must be numeric."")"
$ws.Range("G9").Value = "ValueError(""first_name cannot be blank."")"
$ws.Range("G10").Value = "ValueError(""last_name cannot be blank."")"
$ws.Range("G11").Value = "EmailNotValidError: ""email@pixell-river.com"""

# ...then Method Inputs column (F) for rows 8-11
$ws.Range("F8").Value = "client_number: None`nfirst_name: ""Sion""`nlast_name: ""Kim""`nemail_address: ""skim14@rrc.ca"""
$ws.Range("F9").Value = "client_number: 1999`nfirst_name: "" ""`nlast_name: ""Kim""`nemail_address: ""skim14@rrc.ca"""
$ws.Range("F10").Value = "client_number: 1999`nfirst_name: ""Sion""`nlast_name: "" ""`nemail_address: ""skim14@rrc.ca"""
$ws.Range("F11").Value = "client_number: 1999`nfirst_name: ""Sion""`nlast_name: ""Kim""`nemail_address: ""skim14"""

# Preconditions column (E) for rows 8-11 - reuses "None"
$ws.Range("E8").Value = "None"
$ws.Range("E9").Value = "None"
$ws.Range("E10").Value = "None"
$ws.Range("E11").Value = "None"

# Row 12 - client_number getter
$ws.Range("E12").Value = "Client(1999,`n""Sion"",`n""Kim"",`n""skim14@rrc.ca"")"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "client._Client__client_number`n= 1999"

# Row 13 - first_name getter
$ws.Range("E13").Value = "Client(1999,`n""Sion"",`n""Kim"",`n""skim14@rrc.ca"")"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "client._Client__first_name`n= ""Sion"""

# Row 14 - last_name getter
$ws.Range("E14").Value = "Client(1999,`n""Sion"",`n""Kim"",`n""skim14@rrc.ca"")"
$ws.Range("F14").Value = "None"
$ws.Range("G14").Value = "client._Client__last_name`n= ""Kim"""

# Row 15 - email_address getter
$ws.Range("E15").Value = "Client(1999,`n""Sion"",`n""Kim"",`n""skim14@rrc.ca"")"
$ws.Range("F15").Value = "None"
$ws.Range("G15").Value = "client._Client__email_address`n= ""skim14@rrc.ca"""

# Row 16 - __str__
$ws.Range("E16").Value = "Client(1999,`n""Sion"",`n""Kim"",`n""skim14@rrc.ca"")"
$ws.Range("F16").Value = "None"
$ws.Range("G16").Value = "Kim, Sion [1999] - skim14@rrc.ca"

# Row heights grew to fit the newly entered multi-line content
$ws.Rows(12).RowHeight = 64.5
$ws.Rows(13).RowHeight = 70.5
$ws.Rows(14).RowHeight = 78.4
$ws.Rows(15).RowHeight = 69.4
$ws.Rows(16).RowHeight = 60.75

# Column G (Expected Result) was widened to fit the new, longer text
$ws.Columns(7).ColumnWidth = 29.57

# View state: zoomed out and scrolled down, with a stray selection left on P10
$win = $ws.Application.ActiveWindow
$win.Zoom = 57
$ws.Range("P10").Select() | Out-Null
